$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('T006.207.11.038.00', 'SI2312007', 5, '現貨', 3127, 'cs')
    ,@('T006.207.16.038.00', 'SI2312007', 3, '現貨', 2832, 'cs')
    ,@('T006.407.11.033.00', 'SI2312007', 7, '現貨', 3127, 'cs')
    ,@('T006.407.11.033.02', 'SI2312007', 4, '現貨', 3894, 'cs')
    ,@('T006.407.11.043.00', 'SI2312007', 2, '現貨', 3127, 'cs')
    ,@('T006.407.11.052.00', 'SI2312007', 5, '現貨', 3127, 'cs')
    ,@('T006.407.11.053.00', 'SI2312007', 5, '現貨', 3127, 'cs')
    ,@('T006.407.16.033.00', 'SI2312007', 8, '現貨', 2832, 'cs')
    ,@('T006.407.16.053.00', 'SI2312007', 5, '現貨', 2832, 'cs')
    ,@('T058.009.11.031.00', 'SI2312007', 5, '現貨', 1475, 'cs')
    ,@('T063.617.36.037.00', 'SI2312007', 3, '現貨', 2301, 'cs')
    ,@('T065.430.11.031.00', 'SI2312007', 2, '現貨', 2714, 'cs')
    ,@('T094.210.11.111.00', 'SI2312007', 5, '現貨', 1652, 'cs')
    ,@('T099.407.11.048.00', 'SI2312007', 5, '現貨', 3953, 'cs')
    ,@('T109.407.11.031.00', 'SI2312007', 5, '現貨', 2271.5, 'cs')
    ,@('T116.617.11.047.01', 'SI2312007', 5, '現貨', 1888, 'cs')
    ,@('T120.407.11.091.01', 'SI2312007', 2, '現貨', 3717, 'cs')
    ,@('T120.417.11.041.01', 'SI2312007', 5, '現貨', 2743.5, 'cs')
    ,@('T120.417.11.051.00', 'SI2312007', 3, '現貨', 2743.5, 'cs')
    ,@('T120.417.11.091.01', 'SI2312007', 2, '現貨', 2743.5, 'cs')
    ,@('T120.607.11.041.00', 'SI2312007', 1, '現貨', 5162.5, 'cs')
    ,@('T120.607.11.041.01', 'SI2312007', 1, '現貨', 5162.5, 'cs')
    ,@('T122.210.16.033.00', 'SI2312007', 2, '現貨', 1534, 'cs')
    ,@('T122.423.11.033.00', 'SI2312007', 3, '現貨', 1858.5, 'cs')
    ,@('T125.617.33.051.00', 'SI2312007', 5, '現貨', 2419, 'cs')
    ,@('T126.010.11.013.00', 'SI2312007', 7, '現貨', 2006, 'cs')
    ,@('T126.010.22.013.01', 'SI2312007', 8, '現貨', 2212.5, 'cs')
    ,@('T126.010.36.013.00', 'SI2312007', 5, '現貨', 2006, 'cs')
    ,@('T126.207.11.013.00', 'SI2312007', 7, '現貨', 3481, 'cs')
    ,@('T126.207.22.013.00', 'SI2312007', 5, '現貨', 3805.5, 'cs')
    ,@('T126.207.36.013.00', 'SI2312007', 3, '現貨', 3481, 'cs')
    ,@('T127.407.11.091.01', 'SI2312007', 5, '現貨', 3894, 'cs')
    ,@('T137.407.11.041.00', 'SI2312007', 5, '現貨', 3274.5, 'cs')
    ,@('T137.407.11.091.00', 'SI2312007', 2, '現貨', 3274.5, 'cs')
    ,@('T137.410.11.041.00', 'SI2312007', 5, '現貨', 1740.5, 'cs')
    ,@('T137.410.11.091.00', 'SI2312007', 2, '現貨', 1740.5, 'cs')
    ,@('T137.410.11.091.01', 'SI2312007', 4, '現貨', 1740.5, 'cs')
    ,@('T41.1.183.33', 'SI2312007', 5, '現貨', 2861.5, 'cs')
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
}

$ws.Range("H18").Select()
